$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 245-247 continue the daily series (dates 44319-44321),
# matching the style (date format/border/bold) used by column A in
# the preceding rows.
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)

$dates = @(44319, 44320, 44321)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(245 + $i, 1).Value = $dates[$i]
}

$row245 = @(30, 4, 2, 98, 34, 13, 18, 14, 8, 4, 13, 31, 40, 0, 5, 3, 34, 7, 11, 27, 217, 17, 17, 14, 30, 2, 4, 0, 4, 14, 10, 74, 29, 12, 1, 21, 20, 38, 2, 29, 964, 4, 1, 0, 5, 1, 1, 1, 0)
$row246 = @(29, 4, 2, 87, 30, 14, 12, 13, 8, 4, 13, 30, 34, 0, 4, 3, 23, 5, 9, 23, 219, 17, 17, 14, 24, 3, 3, 0, 5, 14, 10, 70, 30, 9, 2, 21, 18, 30, 2, 30, 897, 4, 1, 1, 3, 1, 1, 1, 0)
$row247 = @(29, 4, 2, 86, 29, 14, 14, 14, 6, 4, 12, 32, 32, 0, 4, 3, 22, 5, 9, 23, 210, 16, 15, 14, 24, 3, 2, 0, 6, 14, 9, 73, 30, 10, 2, 22, 14, 29, 2, 32, 883, 4, 1, 1, 3, 1, 1, 1, 0)

$rows = @($row245, $row246, $row247)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item(245 + $r, $i + 2).Value = $vals[$i]
    }
}
